$wb = $excel.ActiveWorkbook

# Sheet: Restricciones_del_lider
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
$ws.Range("A2").Value = "2.3000000000000003 - x"
$ws.Range("B2").Value = "-3.3000000000000003"
$ws.Range("D2").Value = "0.51"
$ws.Range("A3").Value = "-2.3000000000000003 + x"
$ws.Range("B3").Value = "1.3000000000000003"
$ws.Range("D3").Value = "0.17"

# Sheet: Restricciones_del_follower
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A2").Value = "3.9796747967479664 - 0.8943089430894309y"
$ws.Range("B2").Value = "-4.979674796747966"
$ws.Range("D2").Value = "0.82"
$ws.Range("E2").Value = "-1.1"
$ws.Range("F2").Value = "-1.1"
$ws.Range("A3").Value = "-1.0234999999999999 + 0.22999999999999998y"
$ws.Range("B3").Value = "0.023499999999999854"
$ws.Range("D3").Value = "0.81"
$ws.Range("E3").Value = "7.6"
$ws.Range("F3").Value = "0"

# Sheet: Punto_modificado
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "2.3000000000000003"
$ws.Range("B2").Value = "4.449999999999999"

# Sheet: Vector_bf
$ws = $wb.Worksheets.Item("Vector_bf")
$ws.Range("A2").Value = "2.479033333333337"

# Sheet: Vector_BF
$ws = $wb.Worksheets.Item("Vector_BF")
$ws.Range("A2").Value = "13.386600000000003"
$ws.Range("A3").Value = "-4.5607398373983745"

# Sheet: Vector_Alpha
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 1.23
